# Update countries & provincias Spain
# Refresh the COVID dashboard snapshot: new "last updated" timestamp, a new
# round of per-country figures for the rows whose numbers moved, and a few
# rows where the leaderboard re-sorted by "Casos totales" so the country
# label attached to a given row changed while the row's rank position did
# not.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 00:56"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 3474812
$ws.Range("C4").Value = 60817
$ws.Range("D4").Value = 1544579
$ws.Range("E4").Value = 1792045
$ws.Range("G4").Value = 406
$ws.Range("H4").Value = 138188

# --- Alemania (row 19) ---
$ws.Range("B19").Value = 200436
$ws.Range("C19").Value = 486
$ws.Range("E19").Value = 6197

# --- Colombia (row 22) ---
$ws.Range("B22").Value = 154277
$ws.Range("C22").Value = 3832
$ws.Range("D22").Value = 65809
$ws.Range("E22").Value = 83013
$ws.Range("G22").Value = 148
$ws.Range("H22").Value = 5455

# --- Argentina (row 25) ---
$ws.Range("B25").Value = 103265
$ws.Range("C25").Value = 3099
$ws.Range("E25").Value = 57189
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = 1903

# --- Rows 50-52: leaderboard re-sort swaps the country labels ---
$ws.Range("A50").Value = "Nigeria"
$ws.Range("B50").Value = 33153
$ws.Range("C50").Value = 595
$ws.Range("D50").Value = 13671
$ws.Range("E50").Value = 18738
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 744

$ws.Range("A51").Value = "Rumania"
$ws.Range("B51").Value = 32948
$ws.Range("C51").Value = 413
$ws.Range("D51").Value = 21692
$ws.Range("E51").Value = 9355
$ws.Range("G51").Value = 17
$ws.Range("H51").Value = 1901

$ws.Range("A52").Value = "Suiza"
$ws.Range("B52").Value = 32946
$ws.Range("C52").Value = 63
$ws.Range("D52").Value = 29600
$ws.Range("E52").Value = 1378
$ws.Range("H52").Value = 1968

# --- Guatemala (row 54) ---
$ws.Range("B54").Value = 29742
$ws.Range("C54").Value = 387
$ws.Range("D54").Value = 4321
$ws.Range("E54").Value = 24177
$ws.Range("G54").Value = 25
$ws.Range("H54").Value = 1244

# --- Rows 57-58: Ghana/Azerbaiyan swap ---
$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 24988
$ws.Range("C57").Value = 470
$ws.Range("D57").Value = 21067
$ws.Range("E57").Value = 3782
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 139

$ws.Range("A58").Value = "Azerbaiyan"
$ws.Range("B58").Value = 24570
$ws.Range("C58").Value = 529
$ws.Range("D58").Value = 15640
$ws.Range("E58").Value = 8617
$ws.Range("G58").Value = 7
$ws.Range("H58").Value = 313

# --- Japon (row 59) ---
$ws.Range("B59").Value = 21868
$ws.Range("C59").Value = 366
$ws.Range("D59").Value = 18103
$ws.Range("E59").Value = 2783

# --- Chequia (row 69) ---
$ws.Range("B69").Value = 13238
$ws.Range("C69").Value = 64
$ws.Range("D69").Value = 8373
$ws.Range("E69").Value = 4512

# --- row 78 ---
$ws.Range("B78").Value = 8984
$ws.Range("C78").Value = 3
$ws.Range("E78").Value = 593

# --- row 90 ---
$ws.Range("E90").Value = 5443
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 39

# --- row 93 ---
$ws.Range("B93").Value = 6026
$ws.Range("C93").Value = 84
$ws.Range("D93").Value = 3475
$ws.Range("E93").Value = 2505

# --- Rows 104-105: Somalia/Guinea Ecuatorial swap ---
$ws.Range("A104").Value = "Somalia"
$ws.Range("B104").Value = 3072
$ws.Range("C104").Value = 13
$ws.Range("D104").Value = 1343
$ws.Range("E104").Value = 1636
$ws.Range("H104").Value = 93

$ws.Range("A105").Value = "Guinea Ecuatorial"
$ws.Range("B105").Value = 3071
$ws.Range("D105").Value = 842
$ws.Range("E105").Value = 2178
$ws.Range("H105").Value = 51

# --- Surinam (row 149) ---
$ws.Range("B149").Value = 780
$ws.Range("C149").Value = 39
$ws.Range("D149").Value = 526
$ws.Range("E149").Value = 236

# --- Rows 157-158: Angola/Tanzania swap ---
$ws.Range("A157").Value = "Angola"
$ws.Range("B157").Value = 525
$ws.Range("C157").Value = 42
$ws.Range("D157").Value = 118
$ws.Range("E157").Value = 381
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 26

$ws.Range("A158").Value = "Tanzania"
$ws.Range("B158").Value = 509
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 183
$ws.Range("E158").Value = 305
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 21
